$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price report (2023-03-20, serial 45005) was added for the
# "Femacal de La Calera - Membrillo" sheet. Insert two rows above the most
# recent existing week (row 87) and push everything else down.
$ws.Rows("87:88").Insert()

# New row 87: Especial
$ws.Range("A87").Value = 3
$ws.Range("B87").Value = "Femacal de La Calera"
$ws.Range("C87").Value = "Coquimbo"
$ws.Range("D87").Value = 45005
$ws.Range("E87").Value = 5
$ws.Range("F87").Value = "Fruta"
$ws.Range("G87").Value = 100104
$ws.Range("H87").Value = "Frutos de pepita"
$ws.Range("I87").Value = 100104003
$ws.Range("J87").Value = "Membrillo"
$ws.Range("K87").Value = "Champion"
$ws.Range("L87").Value = "Especial"
$ws.Range("M87").Value = 56
$ws.Range("N87").Value = 18000
$ws.Range("O87").Value = 18000
$ws.Range("P87").Value = 18000
$ws.Range("Q87").Value = "$/caja 18 kilos empedrada"
$ws.Range("R87").Value = "Región de O'Higgins"
$ws.Range("S87").Value = 1000
$ws.Range("T87").Value = 18

# New row 88: Primera
$ws.Range("A88").Value = 3
$ws.Range("B88").Value = "Femacal de La Calera"
$ws.Range("C88").Value = "Coquimbo"
$ws.Range("D88").Value = 45005
$ws.Range("E88").Value = 5
$ws.Range("F88").Value = "Fruta"
$ws.Range("G88").Value = 100104
$ws.Range("H88").Value = "Frutos de pepita"
$ws.Range("I88").Value = 100104003
$ws.Range("J88").Value = "Membrillo"
$ws.Range("K88").Value = "Champion"
$ws.Range("L88").Value = "Primera"
$ws.Range("M88").Value = 67
$ws.Range("N88").Value = 16000
$ws.Range("O88").Value = 16000
$ws.Range("P88").Value = 16000
$ws.Range("Q88").Value = "$/caja 18 kilos empedrada"
$ws.Range("R88").Value = "Región de O'Higgins"
$ws.Range("S88").Value = 889
$ws.Range("T88").Value = 18
